$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 04:52"

# Update country rows whose rank (by Casos totales) shifted and/or whose data was refreshed

# Row 23: Australia
$ws.Cells.Item(23, 1).Value = "Australia"
$ws.Cells.Item(23, 2).Value = 5635
$ws.Cells.Item(23, 3).Value = 85
$ws.Cells.Item(23, 4).Value = 585
$ws.Cells.Item(23, 5).Value = 5016
$ws.Cells.Item(23, 6).Value = 85
$ws.Cells.Item(23, 7).Value = 4
$ws.Cells.Item(23, 8).Value = 34

# Row 96: Honduras
$ws.Cells.Item(96, 1).Value = "Honduras"
$ws.Cells.Item(96, 2).Value = 268
$ws.Cells.Item(96, 3).Value = 4
$ws.Cells.Item(96, 4).Value = 6
$ws.Cells.Item(96, 5).Value = 240
$ws.Cells.Item(96, 6).Value = 10
$ws.Cells.Item(96, 7).Value = 7
$ws.Cells.Item(96, 8).Value = 22

# Row 97: Uzbekistan
$ws.Cells.Item(97, 1).Value = "Uzbekistan"
$ws.Cells.Item(97, 2).Value = 266
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 25
$ws.Cells.Item(97, 5).Value = 239
$ws.Cells.Item(97, 6).Value = 8
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 2

# Row 111: Bolivia
$ws.Cells.Item(111, 1).Value = "Bolivia"
$ws.Cells.Item(111, 2).Value = 157
$ws.Cells.Item(111, 3).Value = 18
$ws.Cells.Item(111, 4).Value = 2
$ws.Cells.Item(111, 5).Value = 145
$ws.Cells.Item(111, 6).Value = 3
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 10

# Row 112: Venezuela
$ws.Cells.Item(112, 1).Value = "Venezuela"
$ws.Cells.Item(112, 2).Value = 155
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 52
$ws.Cells.Item(112, 5).Value = 96
$ws.Cells.Item(112, 6).Value = 6
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 7

# Row 113: Consejo Danes para los Refugiados
$ws.Cells.Item(113, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(113, 2).Value = 154
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 3
$ws.Cells.Item(113, 5).Value = 133
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 18

# Row 114: Martinica
$ws.Cells.Item(114, 1).Value = "Martinica"
$ws.Cells.Item(114, 2).Value = 145
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 27
$ws.Cells.Item(114, 5).Value = 115
$ws.Cells.Item(114, 6).Value = 22
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 3

# Row 115: Niger
$ws.Cells.Item(115, 1).Value = "Niger"
$ws.Cells.Item(115, 2).Value = 144
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 136
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 8

# Row 116: Kirguistan
$ws.Cells.Item(116, 1).Value = "Kirguistan"
$ws.Cells.Item(116, 2).Value = 144
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 9
$ws.Cells.Item(116, 5).Value = 134
$ws.Cells.Item(116, 6).Value = 5
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 1

# Row 124: Paraguay
$ws.Cells.Item(124, 1).Value = "Paraguay"
$ws.Cells.Item(124, 2).Value = 104
$ws.Cells.Item(124, 3).Value = 8
$ws.Cells.Item(124, 4).Value = 12
$ws.Cells.Item(124, 5).Value = 89
$ws.Cells.Item(124, 6).Value = 2
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 3

# Row 125: Trinidad yTobago
$ws.Cells.Item(125, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(125, 2).Value = 103
$ws.Cells.Item(125, 3).Value = 0
$ws.Cells.Item(125, 4).Value = 1
$ws.Cells.Item(125, 5).Value = 96
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 6

# Row 126: Ruanda
$ws.Cells.Item(126, 1).Value = "Ruanda"
$ws.Cells.Item(126, 2).Value = 102
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 102
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 0

# Row 127: Gibraltar
$ws.Cells.Item(127, 1).Value = "Gibraltar"
$ws.Cells.Item(127, 2).Value = 98
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(127, 4).Value = 52
$ws.Cells.Item(127, 5).Value = 46
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

# Row 157: Haiti
$ws.Cells.Item(157, 1).Value = "Haiti"
$ws.Cells.Item(157, 2).Value = 21
$ws.Cells.Item(157, 3).Value = 1
$ws.Cells.Item(157, 4).Value = 1
$ws.Cells.Item(157, 5).Value = 20
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 0

# Row 158: Gabon
$ws.Cells.Item(158, 1).Value = "Gabon"
$ws.Cells.Item(158, 2).Value = 21
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 1
$ws.Cells.Item(158, 5).Value = 19
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 1
